# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the header style (bold, bordered, centered) from an existing header cell (AC1) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2 - 47: every row gets the same team record values
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 77   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 1    # AF -> Ties
}
